# Add season-record columns (Wins / Losses / Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column past the existing "Unnamed: 28" (AC) column.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the other header cells (bold font, borders, centered).
$headerSample = $ws.Range("A1")
$newHeaders = $ws.Range("AD1:AF1")
$headerSample.Copy()
$newHeaders.PasteSpecial(-4122)

# Same season record for every player row (2-44): 68 wins, 94 losses, 0 ties.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 68
    $ws.Cells.Item($r, 31).Value = 94
    $ws.Cells.Item($r, 32).Value = 0
}
